$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 153, shifting existing rows 153:182 down to 154:183.
$ws.Rows.Item(153).Insert()

# Populate the new row 153 with the new weekly record.
$ws.Range("A153").Value = 4
$ws.Range("B153").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C153").Value = "Los Lagos"
$ws.Range("D153").Value = 44543
$ws.Range("E153").Value = 10
$ws.Range("F153").Value = 100112021
$ws.Range("G153").Value = "Ají"
$ws.Range("H153").Value = "Inferno"
$ws.Range("I153").Value = "Primera"
$ws.Range("J153").Value = 35
$ws.Range("K153").Value = 18000
$ws.Range("L153").Value = 18000
$ws.Range("M153").Value = 18000
$ws.Range("N153").Value = "$/caja 12 kilos"
$ws.Range("O153").Value = "Región de Arica y Parinacota"
$ws.Range("P153").Value = 1500
$ws.Range("Q153").Value = 12
$ws.Range("R153").Value = "Hortaliza"
